$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.588.42"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "2.495.13"
$ws.Range("E3").Value = "  -1.00%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.77%  "

$ws.Range("D9").Value = "2.493.19"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("D14").Value = "2.945.74"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").Value = "69.476.98"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.94%  "

$ws.Range("D18").Value = "2.511.18"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "346.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("E23").Value = "  -3.43%  "

$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.49%  "

$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("D28").Value = "2.617.31"
$ws.Range("E28").Value = "  -1.61%  "

$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("E30").Value = "  -2.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "439.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.66%  "

$ws.Range("E33").Value = "  -5.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("E37").Value = "  -3.57%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.56%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.13%  "

$ws.Range("E43").Value = "  -2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.99%  "

$ws.Range("E45").Value = "  -7.28%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.65%  "

$ws.Range("E47").Value = "  -1.70%  "

$ws.Range("E48").Value = "  -3.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0725"
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = "  -1.02%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0922"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.15%  "
